# Merge the split runs ("A" / " " / "slide", "a" / " " / "table",
# "Plus" / " " / "an" / " " / "image") back into single runs per paragraph.
#
# Setting TextRange.Text to the value it already resolves to is a no-op
# for the run structure, so each target is first set to a placeholder
# string and then to the real text -- forcing the host to rewrite the
# paragraph as a single run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1: "Title 1" -> "A slide"
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "__tmp__"
$title.TextFrame.TextRange.Text = "A slide"

# Shape 3: "Content Placeholder 5" (table) -> cell (1,2) "a table"
$tableShape = $s.Shapes.Item(3)
$cell = $tableShape.Table.Cell(1, 2)
$cell.Shape.TextFrame.TextRange.Text = "__tmp__"
$cell.Shape.TextFrame.TextRange.Text = "a table"

# Shape 7: "TextBox 3" -> "Plus an image"
$textBox = $s.Shapes.Item(7)
$textBox.TextFrame.TextRange.Text = "__tmp__"
$textBox.TextFrame.TextRange.Text = "Plus an image"
